$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Neo4j query text for the StudyFilesTab row (note leading space, matches source data)
$newQuery = @"
 MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['COTC022'] and demo.breed in ['Anatolian Shepherd Dog','Saint Bernard'] and diag.disease_term in ['Osteosarcoma'] and diag.primary_disease_site in ['Bone (Appendicular)']
WITH DISTINCT f, s
RETURN 
  coalesce(f.file_name, '') AS ``File Name``,
  coalesce(f.file_type, '') AS ``File Type``,
  coalesce("study", '') AS ``Association``,
  coalesce(f.file_description, '') AS ``Description``,
  coalesce(f.file_format, '') AS ``File Format``,
  coalesce(f.file_size, '') AS ``Size``,
  coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

# Add a new row (row 5) describing a "StudyFilesTab" script, reusing the same
# StatQuery / dbExcel / WebExcel values already used by the other tabs.
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5").Value = $newQuery
$ws.Range("C5").Value = $ws.Range("C4").Text
$ws.Range("D5").Value = $ws.Range("D4").Text
$ws.Range("E5").Value = $ws.Range("E4").Text

# Match the wrap-text style used on the other rows' query/StatQuery columns
$ws.Range("B5:C5").WrapText = $true

# Move the selection/active cell onto the newly added row
$ws.Range("C5").Select()
